# Auto-generated edit script for Uruguay Primera Division workbook
# Applies: rotation of rows 117/119/120 odds data, and update/addition of rows 177-182

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$colmap = @{
  "A" = 1
  "B" = 2
  "C" = 3
  "D" = 4
  "E" = 5
  "F" = 6
  "G" = 7
  "H" = 8
  "I" = 9
  "J" = 10
  "K" = 11
  "L" = 12
  "M" = 13
  "N" = 14
  "O" = 15
  "P" = 16
  "Q" = 17
  "R" = 18
  "S" = 19
  "T" = 20
  "U" = 21
  "V" = 22
  "W" = 23
  "X" = 24
  "Y" = 25
  "Z" = 26
  "AA" = 27
  "AB" = 28
  "AC" = 29
}

function Set-Cell($row, $col, $val) {
    $c = $colmap[$col]
    if ($val -is [string]) {
        $ws.Cells.Item($row, $c).Value = $val
    } else {
        $ws.Cells.Item($row, $c).Value2 = $val
    }
}

function Set-RowData($row, $data) {
    foreach ($pair in $data) {
        Set-Cell $row $pair.Col $pair.Val
    }
}

# --- Update rows 117, 119, 120 (data rotation among three existing matches) ---
$row117 = @(
    @{ Col = "A"; Val = 115 },
    @{ Col = "B"; Val = 7013702 },
    @{ Col = "C"; Val = "Uruguay Primera División" },
    @{ Col = "D"; Val = "Uruguay Clausura" },
    @{ Col = "E"; Val = 45267.70833333334 },
    @{ Col = "F"; Val = "Defensor Sporting" },
    @{ Col = "G"; Val = "Danubio" },
    @{ Col = "H"; Val = 0 },
    @{ Col = "I"; Val = 2 },
    @{ Col = "J"; Val = "A" },
    @{ Col = "K"; Val = 1.8 },
    @{ Col = "L"; Val = 3.6 },
    @{ Col = "M"; Val = 4.2 },
    @{ Col = "N"; Val = 1.8 },
    @{ Col = "O"; Val = 3.6 },
    @{ Col = "P"; Val = 4.2 },
    @{ Col = "Q"; Val = -0.75 },
    @{ Col = "R"; Val = 2.05 },
    @{ Col = "S"; Val = 1.8 },
    @{ Col = "T"; Val = 2.25 },
    @{ Col = "U"; Val = 1.85 },
    @{ Col = "V"; Val = 2 },
    @{ Col = "W"; Val = -1 },
    @{ Col = "X"; Val = -1 },
    @{ Col = "Y"; Val = 3.2 },
    @{ Col = "Z"; Val = -1 },
    @{ Col = "AA"; Val = 0.8 },
    @{ Col = "AB"; Val = -0.5 },
    @{ Col = "AC"; Val = 0.5 }
)

$row119 = @(
    @{ Col = "A"; Val = 117 },
    @{ Col = "B"; Val = 7013409 },
    @{ Col = "C"; Val = "Uruguay Primera División" },
    @{ Col = "D"; Val = "Uruguay Clausura" },
    @{ Col = "E"; Val = 45267.70833333334 },
    @{ Col = "F"; Val = "Nacional De Football" },
    @{ Col = "G"; Val = "Torque" },
    @{ Col = "H"; Val = 1 },
    @{ Col = "I"; Val = 1 },
    @{ Col = "J"; Val = "D" },
    @{ Col = "K"; Val = 1.666 },
    @{ Col = "L"; Val = 3.9 },
    @{ Col = "M"; Val = 4.5 },
    @{ Col = "N"; Val = 1.615 },
    @{ Col = "O"; Val = 4 },
    @{ Col = "P"; Val = 4.75 },
    @{ Col = "Q"; Val = -0.75 },
    @{ Col = "R"; Val = 1.8 },
    @{ Col = "S"; Val = 2.05 },
    @{ Col = "T"; Val = 2.75 },
    @{ Col = "U"; Val = 1.95 },
    @{ Col = "V"; Val = 1.9 },
    @{ Col = "W"; Val = -1 },
    @{ Col = "X"; Val = 3 },
    @{ Col = "Y"; Val = -1 },
    @{ Col = "Z"; Val = -1 },
    @{ Col = "AA"; Val = 1.05 },
    @{ Col = "AB"; Val = -1 },
    @{ Col = "AC"; Val = 0.8999999999999999 }
)

$row120 = @(
    @{ Col = "A"; Val = 118 },
    @{ Col = "B"; Val = 7013886 },
    @{ Col = "C"; Val = "Uruguay Primera División" },
    @{ Col = "D"; Val = "Uruguay Clausura" },
    @{ Col = "E"; Val = 45267.70833333334 },
    @{ Col = "F"; Val = "Racing Club de Montevideo" },
    @{ Col = "G"; Val = "Cerro" },
    @{ Col = "H"; Val = 0 },
    @{ Col = "I"; Val = 1 },
    @{ Col = "J"; Val = "A" },
    @{ Col = "K"; Val = 2.25 },
    @{ Col = "L"; Val = 3.1 },
    @{ Col = "M"; Val = 3.25 },
    @{ Col = "N"; Val = 2.25 },
    @{ Col = "O"; Val = 2.875 },
    @{ Col = "P"; Val = 3.5 },
    @{ Col = "Q"; Val = -0.25 },
    @{ Col = "R"; Val = 1.95 },
    @{ Col = "S"; Val = 1.9 },
    @{ Col = "T"; Val = 2 },
    @{ Col = "U"; Val = 1.925 },
    @{ Col = "V"; Val = 1.925 },
    @{ Col = "W"; Val = -1 },
    @{ Col = "X"; Val = -1 },
    @{ Col = "Y"; Val = 2.5 },
    @{ Col = "Z"; Val = -1 },
    @{ Col = "AA"; Val = 0.8999999999999999 },
    @{ Col = "AB"; Val = -1 },
    @{ Col = "AC"; Val = 0.925 }
)

Set-RowData 117 $row117
Set-RowData 119 $row119
Set-RowData 120 $row120

# --- Update existing rows 177 and 178 with new match results ---
$row177 = @(
    @{ Col = "A"; Val = 175 },
    @{ Col = "B"; Val = 8051184 },
    @{ Col = "C"; Val = "Uruguay Primera División" },
    @{ Col = "D"; Val = "Uruguay Apertura" },
    @{ Col = "E"; Val = 45395.625 },
    @{ Col = "F"; Val = "Club Atletico Progreso" },
    @{ Col = "G"; Val = "Boston River" },
    @{ Col = "H"; Val = 1 },
    @{ Col = "I"; Val = 0 },
    @{ Col = "J"; Val = "H" },
    @{ Col = "K"; Val = 2.375 },
    @{ Col = "L"; Val = 3.2 },
    @{ Col = "M"; Val = 3 },
    @{ Col = "N"; Val = 2.625 },
    @{ Col = "O"; Val = 3.2 },
    @{ Col = "P"; Val = 2.7 },
    @{ Col = "Q"; Val = 0 },
    @{ Col = "R"; Val = 1.925 },
    @{ Col = "S"; Val = 1.925 },
    @{ Col = "T"; Val = 2.25 },
    @{ Col = "U"; Val = 1.95 },
    @{ Col = "V"; Val = 1.9 },
    @{ Col = "W"; Val = 1.625 },
    @{ Col = "X"; Val = -1 },
    @{ Col = "Y"; Val = -1 },
    @{ Col = "Z"; Val = 0.925 },
    @{ Col = "AA"; Val = -1 },
    @{ Col = "AB"; Val = -1 },
    @{ Col = "AC"; Val = 0.8999999999999999 }
)

$row178 = @(
    @{ Col = "A"; Val = 176 },
    @{ Col = "B"; Val = 8051185 },
    @{ Col = "C"; Val = "Uruguay Primera División" },
    @{ Col = "D"; Val = "Uruguay Apertura" },
    @{ Col = "E"; Val = 45395.72916666666 },
    @{ Col = "F"; Val = "Cerro" },
    @{ Col = "G"; Val = "Atletico Fenix Montevideo" },
    @{ Col = "H"; Val = 0 },
    @{ Col = "I"; Val = 0 },
    @{ Col = "J"; Val = "D" },
    @{ Col = "K"; Val = 2.375 },
    @{ Col = "L"; Val = 3 },
    @{ Col = "M"; Val = 3.2 },
    @{ Col = "N"; Val = 2.7 },
    @{ Col = "O"; Val = 3 },
    @{ Col = "P"; Val = 2.8 },
    @{ Col = "Q"; Val = 0 },
    @{ Col = "R"; Val = 1.85 },
    @{ Col = "S"; Val = 2 },
    @{ Col = "T"; Val = 2 },
    @{ Col = "U"; Val = 1.85 },
    @{ Col = "V"; Val = 2 },
    @{ Col = "W"; Val = -1 },
    @{ Col = "X"; Val = 2 },
    @{ Col = "Y"; Val = -1 },
    @{ Col = "Z"; Val = 0 },
    @{ Col = "AA"; Val = -0.0 },
    @{ Col = "AB"; Val = -1 },
    @{ Col = "AC"; Val = 1 }
)

Set-RowData 177 $row177
Set-RowData 178 $row178

# --- Add new rows 179-182: first copy cell formatting (style) for columns A and E
# from an existing fully formatted row, then populate values ---
$formatSourceRow = 176
foreach ($newRow in 179..182) {
    $ws.Cells.Item($formatSourceRow, $colmap["A"]).Copy() | Out-Null
    $ws.Cells.Item($newRow, $colmap["A"]).PasteSpecial(-4122) | Out-Null
    $ws.Cells.Item($formatSourceRow, $colmap["E"]).Copy() | Out-Null
    $ws.Cells.Item($newRow, $colmap["E"]).PasteSpecial(-4122) | Out-Null
}
$excel.CutCopyMode = 0

$row179 = @(
    @{ Col = "A"; Val = 177 },
    @{ Col = "B"; Val = 8051186 },
    @{ Col = "C"; Val = "Uruguay Primera División" },
    @{ Col = "D"; Val = "Uruguay Apertura" },
    @{ Col = "E"; Val = 45395.83333333334 },
    @{ Col = "F"; Val = "Deportivo Maldonado" },
    @{ Col = "G"; Val = "Cerro Largo" },
    @{ Col = "H"; Val = 0 },
    @{ Col = "I"; Val = 1 },
    @{ Col = "J"; Val = "A" },
    @{ Col = "K"; Val = 2.3 },
    @{ Col = "L"; Val = 3 },
    @{ Col = "M"; Val = 3.3 },
    @{ Col = "N"; Val = 2.375 },
    @{ Col = "O"; Val = 3.1 },
    @{ Col = "P"; Val = 3.1 },
    @{ Col = "Q"; Val = -0.25 },
    @{ Col = "R"; Val = 2.05 },
    @{ Col = "S"; Val = 1.8 },
    @{ Col = "T"; Val = 2.25 },
    @{ Col = "U"; Val = 2.1 },
    @{ Col = "V"; Val = 1.775 },
    @{ Col = "W"; Val = -1 },
    @{ Col = "X"; Val = -1 },
    @{ Col = "Y"; Val = 2.1 },
    @{ Col = "Z"; Val = -1 },
    @{ Col = "AA"; Val = 0.8 },
    @{ Col = "AB"; Val = -1 },
    @{ Col = "AC"; Val = 0.7749999999999999 }
)

$row180 = @(
    @{ Col = "A"; Val = 178 },
    @{ Col = "B"; Val = 8051004 },
    @{ Col = "C"; Val = "Uruguay Primera División" },
    @{ Col = "D"; Val = "Uruguay Apertura" },
    @{ Col = "E"; Val = 45396.41666666666 },
    @{ Col = "F"; Val = "Rampla Juniors" },
    @{ Col = "G"; Val = "Racing Club de Montevideo" },
    @{ Col = "H"; Val = 0 },
    @{ Col = "I"; Val = 1 },
    @{ Col = "J"; Val = "A" },
    @{ Col = "K"; Val = 3.2 },
    @{ Col = "L"; Val = 3.3 },
    @{ Col = "M"; Val = 2.2 },
    @{ Col = "N"; Val = 3.8 },
    @{ Col = "O"; Val = 3.6 },
    @{ Col = "P"; Val = 1.85 },
    @{ Col = "Q"; Val = 0.5 },
    @{ Col = "R"; Val = 1.925 },
    @{ Col = "S"; Val = 1.925 },
    @{ Col = "T"; Val = 2.5 },
    @{ Col = "U"; Val = 2.05 },
    @{ Col = "V"; Val = 1.8 },
    @{ Col = "W"; Val = -1 },
    @{ Col = "X"; Val = -1 },
    @{ Col = "Y"; Val = 0.8500000000000001 },
    @{ Col = "Z"; Val = -1 },
    @{ Col = "AA"; Val = 0.925 },
    @{ Col = "AB"; Val = -1 },
    @{ Col = "AC"; Val = 0.8 }
)

$row181 = @(
    @{ Col = "A"; Val = 179 },
    @{ Col = "B"; Val = 8050912 },
    @{ Col = "C"; Val = "Uruguay Primera División" },
    @{ Col = "D"; Val = "Uruguay Apertura" },
    @{ Col = "E"; Val = 45398.75 },
    @{ Col = "F"; Val = "Montevideo Wanderers" },
    @{ Col = "G"; Val = "Liverpool Montevideo" },
    @{ Col = "K"; Val = 3.2 },
    @{ Col = "L"; Val = 3.3 },
    @{ Col = "M"; Val = 2.2 },
    @{ Col = "N"; Val = 3.5 },
    @{ Col = "O"; Val = 3.4 },
    @{ Col = "P"; Val = 2.05 },
    @{ Col = "Q"; Val = 0.25 },
    @{ Col = "R"; Val = 2.025 },
    @{ Col = "S"; Val = 1.825 },
    @{ Col = "T"; Val = 2.25 },
    @{ Col = "U"; Val = 1.925 },
    @{ Col = "V"; Val = 1.925 },
    @{ Col = "W"; Val = 0 },
    @{ Col = "X"; Val = 0 },
    @{ Col = "Y"; Val = 0 },
    @{ Col = "Z"; Val = 0 },
    @{ Col = "AA"; Val = 0 }
)

$row182 = @(
    @{ Col = "A"; Val = 180 },
    @{ Col = "B"; Val = 8050913 },
    @{ Col = "C"; Val = "Uruguay Primera División" },
    @{ Col = "D"; Val = "Uruguay Apertura" },
    @{ Col = "E"; Val = 45398.85416666666 },
    @{ Col = "F"; Val = "Miramar Misiones" },
    @{ Col = "G"; Val = "Nacional De Football" },
    @{ Col = "K"; Val = 5 },
    @{ Col = "L"; Val = 3.5 },
    @{ Col = "M"; Val = 1.727 },
    @{ Col = "N"; Val = 5.75 },
    @{ Col = "O"; Val = 3.75 },
    @{ Col = "P"; Val = 1.6 },
    @{ Col = "Q"; Val = 0.75 },
    @{ Col = "R"; Val = 2.1 },
    @{ Col = "S"; Val = 1.775 },
    @{ Col = "T"; Val = 2.25 },
    @{ Col = "U"; Val = 1.95 },
    @{ Col = "V"; Val = 1.9 },
    @{ Col = "W"; Val = 0 },
    @{ Col = "X"; Val = 0 },
    @{ Col = "Y"; Val = 0 },
    @{ Col = "Z"; Val = 0 },
    @{ Col = "AA"; Val = 0 }
)

Set-RowData 179 $row179
Set-RowData 180 $row180
Set-RowData 181 $row181
Set-RowData 182 $row182

